$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.486.59"
$ws.Range("E2").Value = "  +1.28%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.16"
$ws.Range("E3").Value = "  +2.21%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("B5").Value = "USDC"
$ws.Range("C5").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("E5").Value = "  +0.08%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.85"
$ws.Range("E6").Value = "  +1.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3767"
$ws.Range("E7").Value = "  -0.35%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.44"
$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3645"
$ws.Range("E9").Value = "  +0.90%  "

$ws.Range("E10").Value = "  +0.40%  "

$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.98"
$ws.Range("E13").Value = "  +1.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.641"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001280"
$ws.Range("E15").Value = "  +2.59%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.381"
$ws.Range("E16").Value = "  -0.12%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.637.97"
$ws.Range("E17").Value = "  +2.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.77"
$ws.Range("E18").Value = "  +1.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06969"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.21"
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.545"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.483.88"
$ws.Range("E23").Value = "  +1.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.82"
$ws.Range("E24").Value = "  -1.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.087"
$ws.Range("E25").Value = "  +3.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.421"
$ws.Range("E26").Value = "  +1.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.27"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.27"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.347"
$ws.Range("E29").Value = "  +1.95%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.66"
$ws.Range("E30").Value = "  +1.37%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.358"
$ws.Range("E31").Value = "  -1.22%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.819.02"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.786"
$ws.Range("E33").Value = "  -0.93%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9652"
$ws.Range("E34").Value = "  -0.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02821"
$ws.Range("E35").Value = "  +3.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.35"
$ws.Range("E36").Value = "  +0.44%  "

$ws.Range("E37").Value = "  -2.26%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2537"
$ws.Range("E38").Value = "  +1.17%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.182"
$ws.Range("E39").Value = "  +1.07%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08862"
$ws.Range("E40").Value = "  +0.70%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.381"
$ws.Range("E41").Value = "  +1.45%  "

$ws.Range("E42").Value = "  +0.18%  "

$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.52"
$ws.Range("E43").Value = "  +0.11%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.23"
$ws.Range("E44").Value = "  +5.45%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6551"
$ws.Range("E45").Value = "  +0.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.340"
$ws.Range("E46").Value = "  +1.45%  "

$ws.Range("E47").Value = "  +0.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.029"
$ws.Range("E48").Value = "  +0.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07974"
$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.71"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.209"
$ws.Range("E51").Value = "  +0.41%  "
